$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The periodic "Opdateret d. <dato>" refresh: the sheet (and the workbook-
# level defined name that mirrors it) is renamed to the new update date.
# Excel keeps the defined name's sheet-qualified reference in sync with the
# sheet name automatically.
$ws.Name = "Opdateret d. 05-12-2025"
